$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "67.812.18"
$ws.Cells.Item(2, 5).Value = "  +3.42%  "
# Row 3
$ws.Cells.Item(3, 4).Value = "3.593.54"
$ws.Cells.Item(3, 5).Value = "  +3.17%  "
# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.53%  "
# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "199.72"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +9.72%  "
# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "569.15"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  +1.23%  "
# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.613"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = "  +2.02%  "
# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.00%  "
# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.676"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = "  +4.08%  "
# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "59.93"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = "  +15.88%  "
# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.147"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  +3.57%  "
# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0000280"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +11.56%  "
# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "10.25"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = "  +7.35%  "
# Row 14
$ws.Cells.Item(14, 4).Value = "4.186.23"
$ws.Cells.Item(14, 5).Value = "  +3.88%  "
# Row 15
$ws.Cells.Item(15, 4).Value = "3.612.80"
$ws.Cells.Item(15, 5).Value = "  +4.11%  "
# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.126"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value = "  +1.45%  "
# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "19.14"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).Value = "  +7.74%  "
# Row 18
$ws.Cells.Item(18, 4).Value = "67.740.23"
$ws.Cells.Item(18, 5).Value = "  +3.69%  "
# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.21"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +4.13%  "
# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "1.06"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "  +2.37%  "
# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "402.82"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = "  +6.23%  "
# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "12.85"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = "  +20.29%  "
# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.18"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "  +2.21%  "
# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "84.82"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = "  +1.35%  "
# Row 25
$ws.Cells.Item(25, 2).Value = "ImmutableX"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.89"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "  +2.00%  "
# Row 26
$ws.Cells.Item(26, 2).Value = "Toncoin"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.90"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = "  +12.59%  "
# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "12.42"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = "  +3.75%  "
# Row 28
$ws.Cells.Item(28, 5).Value = "  +2.25%  "
# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.20"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "  +6.41%  "
# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "7.66"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  +5.83%  "
# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "31.44"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +3.12%  "
# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "676.15"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = "  +11.79%  "
# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "12.11"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = "  +1.32%  "
# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.113"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  +2.26%  "
# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "63.37"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = "  +1.99%  "
# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "41.34"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "  +1.39%  "
# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.409"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  +4.20%  "
# Row 38
$ws.Cells.Item(38, 5).Value = "  -0.01%  "
# Row 39
$ws.Cells.Item(39, 2).Value = "ThetaToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.20"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = "  +16.25%  "
# Row 40
$ws.Cells.Item(40, 2).Value = "PEPE"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(40, 4).Value = "0.0₃0757"
$ws.Cells.Item(40, 5).Value = "  +5.81%  "
# Row 41
$ws.Cells.Item(41, 4).Value = "3.182.75"
$ws.Cells.Item(41, 5).Value = "  +7.41%  "
# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.133"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  +3.12%  "
# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.21%  "
# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.68"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = "  +9.46%  "
# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.81"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +15.39%  "
# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.78"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "  +20.51%  "
# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0410"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = "  +3.53%  "
# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.131"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = "  +2.70%  "
# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.10"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -1.63%  "
# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "8.63"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = "  +5.83%  "
# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "139.10"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = "  +2.18%  "
